# Auto-generated edit script: updates cryptos price/volume table
# to match the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force text interpretation so numeric-looking strings (e.g. "54.70",
    # "1.002") keep trailing zeros / do not get coerced to a Double,
    # then restore the default "Normal" style so no stray number format
    # is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "25.848.04"
Set-TextValue $ws.Range("E2") "  -1.23%  "
Set-TextValue $ws.Range("D3") "1.634.38"
Set-TextValue $ws.Range("E3") "  -1.34%  "
Set-TextValue $ws.Range("E4") "  -0.35%  "
Set-TextValue $ws.Range("D5") "214.95"
Set-TextValue $ws.Range("E5") "  -0.55%  "
Set-TextValue $ws.Range("D6") "0.5014"
Set-TextValue $ws.Range("E6") "  -1.70%  "
Set-TextValue $ws.Range("E7") "  -0.36%  "
Set-TextValue $ws.Range("D8") "0.2568"
Set-TextValue $ws.Range("E8") "  -0.68%  "
Set-TextValue $ws.Range("D9") "0.06399"
Set-TextValue $ws.Range("E9") "  -0.18%  "
Set-TextValue $ws.Range("D10") "19.59"
Set-TextValue $ws.Range("E10") "  -2.00%  "
Set-TextValue $ws.Range("D11") "0.07698"
Set-TextValue $ws.Range("E11") "  -1.30%  "
Set-TextValue $ws.Range("B12") "WrappedEther"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D12") "1.638.66"
Set-TextValue $ws.Range("E12") "  -1.01%  "
Set-TextValue $ws.Range("B13") "Polkadot"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D13") "4.236"
Set-TextValue $ws.Range("E13") "  -1.34%  "
Set-TextValue $ws.Range("D14") "1.860.61"
Set-TextValue $ws.Range("D15") "0.5426"
Set-TextValue $ws.Range("E15") "  -2.23%  "
Set-TextValue $ws.Range("D16") "0.0₅7920"
Set-TextValue $ws.Range("E16") "  -1.32%  "
Set-TextValue $ws.Range("D17") "63.41"
Set-TextValue $ws.Range("E17") "  -1.15%  "
Set-TextValue $ws.Range("D18") "25.860.35"
Set-TextValue $ws.Range("E18") "  -1.31%  "
Set-TextValue $ws.Range("D19") "1.002"
Set-TextValue $ws.Range("E19") "  -0.38%  "
Set-TextValue $ws.Range("D20") "202.78"
Set-TextValue $ws.Range("E20") "  -3.30%  "
Set-TextValue $ws.Range("D21") "4.317"
Set-TextValue $ws.Range("E21") "  -2.24%  "
Set-TextValue $ws.Range("D22") "9.930"
Set-TextValue $ws.Range("E22") "  -1.48%  "
Set-TextValue $ws.Range("D23") "5.971"
Set-TextValue $ws.Range("E23") "  -0.87%  "
Set-TextValue $ws.Range("E24") "  -0.25%  "
Set-TextValue $ws.Range("E25") "  +11.08%  "
Set-TextValue $ws.Range("D26") "140.73"
Set-TextValue $ws.Range("E26") "  -2.27%  "
Set-TextValue $ws.Range("D27") "0.1141"
Set-TextValue $ws.Range("E27") "  -2.40%  "
Set-TextValue $ws.Range("D28") "15.67"
Set-TextValue $ws.Range("E28") "  -0.67%  "
Set-TextValue $ws.Range("D29") "6.697"
Set-TextValue $ws.Range("E29") "  -4.06%  "
Set-TextValue $ws.Range("D30") "1.238"
Set-TextValue $ws.Range("E30") "  -0.85%  "
Set-TextValue $ws.Range("D31") "0.04979"
Set-TextValue $ws.Range("E31") "  -2.73%  "
Set-TextValue $ws.Range("D32") "3.255"
Set-TextValue $ws.Range("E32") "  -2.76%  "
Set-TextValue $ws.Range("D33") "3.174"
Set-TextValue $ws.Range("E33") "  -1.58%  "
Set-TextValue $ws.Range("D34") "1.537"
Set-TextValue $ws.Range("E34") "  -1.65%  "
Set-TextValue $ws.Range("D35") "2.364"
Set-TextValue $ws.Range("E35") "  -0.43%  "
Set-TextValue $ws.Range("D36") "1.167.13"
Set-TextValue $ws.Range("E36") "  +0.45%  "
Set-TextValue $ws.Range("D37") "0.8919"
Set-TextValue $ws.Range("E37") "  -3.93%  "
Set-TextValue $ws.Range("D38") "2.613"
Set-TextValue $ws.Range("E38") "  -5.00%  "
Set-TextValue $ws.Range("D39") "0.5609"
Set-TextValue $ws.Range("E39") "  -1.63%  "
Set-TextValue $ws.Range("D40") "0.01556"
Set-TextValue $ws.Range("E40") "  -2.33%  "
Set-TextValue $ws.Range("D41") "2.550"
Set-TextValue $ws.Range("E41") "  -0.52%  "
Set-TextValue $ws.Range("D43") "5.671"
Set-TextValue $ws.Range("E43") "  +0.35%  "
Set-TextValue $ws.Range("D44") "0.8073"
Set-TextValue $ws.Range("E44") "  -3.61%  "
Set-TextValue $ws.Range("D45") "99.21"
Set-TextValue $ws.Range("E45") "  -1.41%  "
Set-TextValue $ws.Range("D46") "1.772.71"
Set-TextValue $ws.Range("E46") "  -1.28%  "
Set-TextValue $ws.Range("D47") "0.0₈114"
Set-TextValue $ws.Range("E47") "  -0.97%  "
Set-TextValue $ws.Range("E48") "  -0.70%  "
Set-TextValue $ws.Range("D49") "1.003"
Set-TextValue $ws.Range("D50") "54.70"
Set-TextValue $ws.Range("E50") "  -2.07%  "
Set-TextValue $ws.Range("D51") "0.05082"
Set-TextValue $ws.Range("E51") "  +0.59%  "
